$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (column C) date value for rows 2-7 from 45212 to 45221
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
